$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Imports")

# Update parent/IDs term lists (fixing parents etc.)
$ws.Range("D9").Value = "whole plant [PO:0000003]"
$ws.Range("D19").Value = "process [BFO:0000015]; object aggregate [BFO:0000027]; role [BFO:0000023]; disposition [BFO:0000016]; object [BFO:0000030]; process profile [BFO:0000144]; site [BFO:0000029]; occurrent [BFO:0000003]; process boundary [BFO:0000035]"
$ws.Range("D20").Value = "is about [IAO:0000136]; data item [IAO:0000027]; report [IAO:0000088]; plan specification [IAO:0000104]; material information bearer [IAO:0000178]; document [IAO:0000310]"
$ws.Range("D3").Value = "Credible interval [OBCS:0000071]; Mortality ratio [OBCS:0000150]; model [OBCS:0000035]"

# Leave selection on D4, matching the saved cursor position
$ws.Range("D4").Select()
